$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 16450
$ws.Range("I13").Value = 2900
$ws.Range("J13").Value = 30000
$ws.Range("K13").Value = 2900
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = -2731
$ws.Range("N13").Value = -30338

$ws.Range("H32").Value = 444.75
$ws.Range("I32").Value = 450.25
$ws.Range("J32").Value = 442.91666
$ws.Range("K32").Value = 450.25
$ws.Range("L32").Value = 442.91666
$ws.Range("M32").Value = -124.25
$ws.Range("N32").Value = -1094.91666

$ws.Range("H38").Value = 2264.5
$ws.Range("I38").Value = 981.875
$ws.Range("J38").Value = 3974.6667
$ws.Range("K38").Value = 2945.625
$ws.Range("L38").Value = 11924.0001
$ws.Range("M38").Value = -2573.625
$ws.Range("N38").Value = -12668.0001

$ws.Range("H40").Value = 6331090.5
$ws.Range("I40").Value = 2003.5146
$ws.Range("J40").Value = 45456356
$ws.Range("K40").Value = 2003.5146
$ws.Range("L40").Value = 45456356
$ws.Range("M40").Value = -1828.5146
$ws.Range("N40").Value = -45456706

$ws.Range("H43").Value = 1594.2273
$ws.Range("I43").Value = 1549.1
$ws.Range("J43").Value = 1631.8334
$ws.Range("K43").Value = 1549.1
$ws.Range("L43").Value = 1631.8334
$ws.Range("M43").Value = -1480.1
$ws.Range("N43").Value = -1769.8334

$ws.Range("H51").Value = 2497.5
$ws.Range("I51").Value = 2496.6667
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 2496.6667
$ws.Range("L51").Value = 2500
$ws.Range("M51").Value = -2012.6667
$ws.Range("N51").Value = -3468

$ws.Range("H80").Value = 5383
$ws.Range("I80").Value = 432.83334
$ws.Range("J80").Value = 12808.25
$ws.Range("K80").Value = 1298.50002
$ws.Range("L80").Value = 38424.75
$ws.Range("M80").Value = -300.5000199999999
$ws.Range("N80").Value = -40420.75

$ws.Range("H83").Value = 5383
$ws.Range("I83").Value = 432.83334
$ws.Range("J83").Value = 12808.25
$ws.Range("K83").Value = 3895.50006
$ws.Range("L83").Value = 115274.25
$ws.Range("M83").Value = 1096.49994
$ws.Range("N83").Value = -125258.25

$ws.Range("H116").Value = 7459.227
$ws.Range("I116").Value = 10313.25
$ws.Range("J116").Value = 4034.4
$ws.Range("K116").Value = 10313.25
$ws.Range("L116").Value = 4034.4
$ws.Range("M116").Value = -6871.25
$ws.Range("N116").Value = -10918.4

$ws.Range("H132").Value = 1620.5526
$ws.Range("I132").Value = 1478.1892
$ws.Range("J132").Value = 6888
$ws.Range("K132").Value = 4434.5676
$ws.Range("L132").Value = 20664
$ws.Range("M132").Value = -1904.5676
$ws.Range("N132").Value = -25724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 929.7742
$ws.Range("I2").Value = 791.36
$ws.Range("J2").Value = 1506.5
$ws.Range("K2").Value = 791.36
$ws.Range("L2").Value = 1506.5
$ws.Range("M2").Value = -678.36
$ws.Range("N2").Value = -1732.5

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = 0

$ws.Range("H116").Value = 929.7742
$ws.Range("I116").Value = 791.36
$ws.Range("J116").Value = 1506.5
$ws.Range("K116").Value = 791.36
$ws.Range("L116").Value = 1506.5
$ws.Range("M116").Value = 1502.64
$ws.Range("N116").Value = -6094.5

$ws.Range("H134").Value = 0
$ws.Range("I134").ClearContents()
$ws.Range("J134").Value = 0
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 929.7742
$ws.Range("I3").Value = 791.36
$ws.Range("J3").Value = 1506.5
$ws.Range("K3").Value = 791.36
$ws.Range("L3").Value = 1506.5
$ws.Range("M3").Value = -677.36
$ws.Range("N3").Value = -1734.5

$ws.Range("H22").Value = 283.41666
$ws.Range("I22").Value = 272.1
$ws.Range("J22").Value = 340
$ws.Range("K22").Value = 272.1
$ws.Range("L22").Value = 340
$ws.Range("M22").Value = -99.10000000000002
$ws.Range("N22").Value = -686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4345.478
$ws.Range("I31").Value = 2010.079
$ws.Range("J31").Value = 7208.2256
$ws.Range("K31").Value = 2010.079
$ws.Range("L31").Value = 7208.2256
$ws.Range("M31").Value = -1715.079
$ws.Range("N31").Value = -7798.2256

$ws.Range("H34").Value = 4345.478
$ws.Range("I34").Value = 2010.079
$ws.Range("J34").Value = 7208.2256
$ws.Range("K34").Value = 2010.079
$ws.Range("L34").Value = 7208.2256
$ws.Range("M34").Value = -1808.079
$ws.Range("N34").Value = -7612.2256

$ws.Range("H111").Value = 23950
$ws.Range("I111").ClearContents()
$ws.Range("J111").Value = 23950
$ws.Range("K111").ClearContents()
$ws.Range("L111").Value = 23950
$ws.Range("N111").Value = -32130

$ws.Range("H115").Value = 25300
$ws.Range("I115").ClearContents()
$ws.Range("J115").Value = 25300
$ws.Range("K115").ClearContents()
$ws.Range("L115").Value = 25300
$ws.Range("N115").Value = -27650

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 214.75
$ws.Range("I10").Value = 186.33333
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 558.99999
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = -419.99999
$ws.Range("N10").Value = -1178

$ws.Range("H20").Value = 5422.25
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 6125.4287
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 18376.2861
$ws.Range("M20").Value = -1273
$ws.Range("N20").Value = -18830.2861

$ws.Range("H23").Value = 5263273
$ws.Range("I23").Value = 16666751
$ws.Range("J23").Value = 129.3077
$ws.Range("K23").Value = 50000253
$ws.Range("L23").Value = 387.9231
$ws.Range("M23").Value = -50000018
$ws.Range("N23").Value = -857.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 50290
$ws.Range("I32").ClearContents()
$ws.Range("J32").Value = 50290
$ws.Range("K32").ClearContents()
$ws.Range("L32").Value = 50290
$ws.Range("N32").Value = -50882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 660.2
$ws.Range("I9").Value = 345.5
$ws.Range("J9").Value = 870
$ws.Range("K9").Value = 345.5
$ws.Range("L9").Value = 870
$ws.Range("M9").Value = -121.5
$ws.Range("N9").Value = -1318

$ws.Range("H22").Value = 1789.9032
$ws.Range("I22").Value = 448.375
$ws.Range("J22").Value = 2256.5217
$ws.Range("K22").Value = 448.375
$ws.Range("L22").Value = 2256.5217
$ws.Range("M22").Value = -153.375
$ws.Range("N22").Value = -2846.5217

$ws.Range("H27").Value = 1789.9032
$ws.Range("I27").Value = 448.375
$ws.Range("J27").Value = 2256.5217
$ws.Range("K27").Value = 448.375
$ws.Range("L27").Value = 2256.5217
$ws.Range("M27").Value = -341.375
$ws.Range("N27").Value = -2470.5217

$ws.Range("H43").Value = 0
$ws.Range("I43").ClearContents()
$ws.Range("J43").Value = 0
$ws.Range("K43").ClearContents()
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0

$ws.Range("H55").Value = 267.33334
$ws.Range("I55").Value = 267.33334
$ws.Range("J55").ClearContents()
$ws.Range("K55").Value = 267.33334
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = -94.33334000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1000000000
$ws.Range("I12").Value = 1000000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1000000000
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -999999858

$ws.Range("H21").Value = 0
$ws.Range("I21").ClearContents()
$ws.Range("J21").Value = 0
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0

$ws.Range("H25").Value = 0
$ws.Range("I25").ClearContents()
$ws.Range("J25").Value = 0
$ws.Range("K25").ClearContents()
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0

$ws.Range("H35").Value = 0
$ws.Range("I35").ClearContents()
$ws.Range("J35").Value = 0
$ws.Range("K35").ClearContents()
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0

$ws.Range("H37").Value = 15000
$ws.Range("I37").ClearContents()
$ws.Range("J37").Value = 15000
$ws.Range("K37").ClearContents()
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15406

$ws.Range("H108").Value = 120000
$ws.Range("I108").ClearContents()
$ws.Range("J108").Value = 120000
$ws.Range("K108").ClearContents()
$ws.Range("L108").Value = 120000
$ws.Range("N108").Value = -127680
